# Update "Generate Report for Handback" timestamps.
#
# The workbook has three worksheets:
#   - Overview : column G "Latest HO Xliff Generate Date"
#   - zh-cn    : column H "Correspond Handoff Datetime", column K "Correspond Handback DateTime"
#   - de-de    : column K "Correspond Handback DateTime"
#
# Row 3 on each sheet corresponds to file 4f570b87-b645-4807-a2d7-5d3068cdcf33,
# whose handback-report timestamps were regenerated for this commit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (G3)
$wsOverview.Range("G3").Value = "2016-08-19 16:49:27"

# zh-cn sheet: Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsZhCn.Range("H3").Value = "2016-08-19 16:49:23"
$wsZhCn.Range("K3").Value = "2016-08-19 16:49:40"

# de-de sheet: Correspond Handback DateTime (K3)
$wsDeDe.Range("K3").Value = "2016-08-19 16:49:47"
